$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a D-column (Price) value to remain plain text, matching the
# workbook's existing convention (inline/shared strings, not numbers), even
# when the text looks numeric (e.g. "0.620", "100.40") so trailing zeros and
# formatting survive. Uses a leading apostrophe to defeat Excel's automatic
# number conversion, then clears the resulting quote-prefix style so the cell
# keeps its original (default) style.
function Set-TextCell($Cell, $Text) {
    $ws.Range($Cell).Value = "'" + $Text
    $ws.Range($Cell).Style = "Normal"
}

Set-TextCell 'D2' '42.357.03'
$ws.Range('E2').Value = '  -2.78%  '

Set-TextCell 'D3' '2.218.63'
$ws.Range('E3').Value = '  -2.28%  '

$ws.Range('E4').Value = '  +0.26%  '

Set-TextCell 'D5' '109.43'
$ws.Range('E5').Value = '  -7.91%  '

Set-TextCell 'D6' '287.19'
$ws.Range('E6').Value = '  +6.96%  '

Set-TextCell 'D7' '0.620'
$ws.Range('E7').Value = '  -3.45%  '

$ws.Range('E8').Value = '  -0.44%  '

$ws.Range('E9').Value = '  -4.82%  '

$ws.Range('E10').Value = '  -9.43%  '

$ws.Range('E11').Value = '  -4.14%  '

Set-TextCell 'D12' '54.14'
$ws.Range('E12').Value = '  +0.02%  '

Set-TextCell 'D13' '8.51'
$ws.Range('E13').Value = '  -9.23%  '

Set-TextCell 'D14' '0.995'
$ws.Range('E14').Value = '  +7.98%  '

$ws.Range('E15').Value = '  -3.06%  '

Set-TextCell 'D16' '14.76'
$ws.Range('E16').Value = '  -7.10%  '

Set-TextCell 'D17' '2.546.44'
$ws.Range('E17').Value = '  -2.57%  '

Set-TextCell 'D18' '2.220.65'
$ws.Range('E18').Value = '  -2.26%  '

Set-TextCell 'D19' '42.188.62'
$ws.Range('E19').Value = '  -3.29%  '

Set-TextCell 'D22' '72.51'
$ws.Range('E22').Value = '  +0.02%  '

Set-TextCell 'D23' '3.33'
$ws.Range('E23').Value = '  +10.70%  '

$ws.Range('E24').Value = '  -1.94%  '

Set-TextCell 'D25' '228.64'
$ws.Range('E25').Value = '  -2.53%  '

Set-TextCell 'D26' '8.79'
$ws.Range('E26').Value = '  -8.42%  '

$ws.Range('E27').Value = '  -1.85%  '

Set-TextCell 'D28' '11.27'
$ws.Range('E28').Value = '  -8.47%  '

$ws.Range('E29').Value = '  -2.51%  '

Set-TextCell 'D30' '172.72'
$ws.Range('E30').Value = '  -1.07%  '

Set-TextCell 'D31' '3.14'
$ws.Range('E31').Value = '  -6.46%  '

Set-TextCell 'D32' '36.58'
$ws.Range('E32').Value = '  -12.28%  '

Set-TextCell 'D33' '20.67'
$ws.Range('E33').Value = '  -3.92%  '

Set-TextCell 'D34' '0.0868'
$ws.Range('E34').Value = '  -5.50%  '

Set-TextCell 'D35' '5.54'
$ws.Range('E35').Value = '  -3.05%  '

$ws.Range('E36').Value = '  +4.89%  '

$ws.Range('E37').Value = '  -4.18%  '

$ws.Range('E38').Value = '  -4.84%  '

Set-TextCell 'D39' '0.0363'
$ws.Range('E39').Value = '  -4.07%  '

Set-TextCell 'D40' '0.102'
$ws.Range('E40').Value = '  -6.36%  '

Set-TextCell 'D41' '72.86'
$ws.Range('E41').Value = '  +1.26%  '

$ws.Range('E42').Value = '  -8.27%  '

$ws.Range('E43').Value = '  -5.10%  '

Set-TextCell 'D45' '12.08'
$ws.Range('E45').Value = '  -11.99%  '

$ws.Range('E46').Value = '  -6.97%  '

Set-TextCell 'D47' '5.27'
$ws.Range('E47').Value = '  -7.96%  '

Set-TextCell 'D48' '1.69'
$ws.Range('E48').Value = '  +5.90%  '

$ws.Range('E49').Value = '  -1.96%  '

$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextCell 'D20' '0.0000104'
$ws.Range('E20').Value = '  -5.17%  '

$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextCell 'D21' '7.04'
$ws.Range('E21').Value = '  +2.18%  '

$ws.Range('B50').Value = 'Aave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextCell 'D50' '100.40'
$ws.Range('E50').Value = '  -2.35%  '

$ws.Range('B51').Value = 'FraxShare'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextCell 'D51' '8.34'
$ws.Range('E51').Value = '  -2.77%  '
